$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 1: labels/offsets header row ---
$ws.Range("A1").Value = "SELLER ORIGIN"
$ws.Range("B1").Formula = "=__METADATA__!B3"

$ws.Range("J1").Value = -28
$ws.Range("K1").Value = -21
$ws.Range("L1").Value = -14
$ws.Range("M1").Value = -7
$ws.Range("N1").Value = 0
$ws.Range("O1").Value = -4
$ws.Range("P1").Value = -3
$ws.Range("Q1").Value = -2
$ws.Range("R1").Value = -1
$ws.Range("S1").Value = -2
$ws.Range("T1").Value = -1
$ws.Range("U1").Value = 0

# --- Row 2: dates derived from B2 (moved from A2) and the row-1 offsets ---
$ws.Range("A2").Value = "DATE"
$ws.Range("B2").Formula = "=__METADATA__!B1"

$ws.Range("J2").Formula = "=B2 + J1"
$ws.Range("K2").Formula = "=B2 + K1"
$ws.Range("L2").Formula = "=B2 + L1"
$ws.Range("M2").Formula = "=B2 + M1"
$ws.Range("N2").Formula = "=B2 + N1"
$ws.Range("O2").Formula = "=EOMONTH(B2, O1)"
$ws.Range("P2").Formula = "=EOMONTH(B2, P1)"
$ws.Range("Q2").Formula = "=EOMONTH(B2, Q1)"
$ws.Range("R2").Formula = "=EOMONTH(B2, R1)"
$ws.Range("S2").Formula = "=EOMONTH(B2, S1)"
$ws.Range("T2").Formula = "=EOMONTH(B2, T1)"
$ws.Range("U2").Formula = "=B2 + U1"

# --- Row 3: region selector plus new ISO week / month / quarter formulas ---
$ws.Range("A3").Value = "REGION"
$ws.Range("B3").Value = "ALL"

$ws.Range("J3").Formula = "=ISOWEEKNUM(J2 + 1)"
$ws.Range("K3").Formula = "=ISOWEEKNUM(K2 + 1)"
$ws.Range("L3").Formula = "=ISOWEEKNUM(L2 + 1)"
$ws.Range("M3").Formula = "=ISOWEEKNUM(M2 + 1)"
$ws.Range("N3").Formula = "=ISOWEEKNUM(N2 + 1)"
$ws.Range("O3").Formula = "=MONTH(O2)"
$ws.Range("P3").Formula = "=MONTH(P2)"
$ws.Range("Q3").Formula = "=MONTH(Q2)"
$ws.Range("R3").Formula = "=MONTH(R2)"
$ws.Range("S3").Formula = "=ROUNDUP(MONTH(S2)/3,0)"
$ws.Range("T3").Formula = "=ROUNDUP(MONTH(T2)/3,0)"
$ws.Range("U3").Formula = "=ISOWEEKNUM(U2 + 1)"

# --- Row 4: marketplace selector ---
$ws.Range("A4").Value = "MARKETPLACE"
$ws.Range("B4").Value = "ALL"

# --- Row 5: team selector ---
$ws.Range("A5").Value = "TEAM"
$ws.Range("B5").Value = "ALL"
